$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (2025-10-26) right after the last existing row.
# Force the date cell to text first so the "MM/DD/YYYY"-looking string is
# kept as a literal string (matching the rest of the column) instead of
# being auto-converted into a date serial value, then restore the cell's
# style to the default so no stray formatting is left behind.
$dateCell = $ws.Range("A70")
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/26/2025"
$dateCell.Style = "Normal"

$ws.Range("B70").Value = 11969.3
